# Update "想去人数" (F column) values on sheets "展览" and "全部类型"
# per the target diff.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> hashtable of row -> new value
$updates = @{
    "展览" = @{
        3  = 99
        4  = 76
        6  = 160
        7  = 7580
        8  = 83
        13 = 428
        14 = 159
        17 = 56
        19 = 23
        20 = 5432
        21 = 138
        22 = 193
        23 = 856
        24 = 223
        25 = 285
    }
    "全部类型" = @{
        3  = 99
        4  = 76
        6  = 160
        7  = 7580
        8  = 83
        13 = 428
        14 = 159
        17 = 56
        19 = 23
        21 = 5432
        23 = 138
        24 = 193
        25 = 856
        26 = 223
        27 = 285
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
